$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.075.09"

$ws.Range("D3").Value = "1.664.27"

$ws.Range("E4").Value = "  -0.83%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5152"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.32%  "

$ws.Range("E7").Value = "  -0.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2628"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06199"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07502"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.30%  "

$ws.Range("D12").Value = "1.668.06"
$ws.Range("E12").Value = "  -1.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.402"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5564"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000007897"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").Value = "26.084.93"
$ws.Range("E17").Value = "  -1.29%  "

$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.768"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "185.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.121"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("E23").Value = "  -0.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "146.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1236"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.506"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06271"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.354"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.00%  "

$ws.Range("E30").Value = "  -4.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.464"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.403"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.608"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9913"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.408"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6006"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.59%  "

$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.087"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01605"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.28%  "

$ws.Range("D40").Value = "1.074.38"
$ws.Range("E40").Value = "  -3.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8557"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.98%  "

$ws.Range("E42").Value = "  -1.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.51%  "

$ws.Range("D44").Value = "1.812.35"
$ws.Range("E44").Value = "  -1.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000108"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05248"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.934"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4248"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.883"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
